$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 236.5
$ws.Cells.Item(2, 9).Value = 79.75
$ws.Cells.Item(2, 10).Value = 550
$ws.Cells.Item(2, 11).Value = 79.75
$ws.Cells.Item(2, 12).Value = 550
$ws.Cells.Item(2, 13).Value = 33.25
$ws.Cells.Item(2, 14).Value = -776
$ws.Cells.Item(18, 8).Value = 558.3333
$ws.Cells.Item(18, 9).Value = 562.5
$ws.Cells.Item(18, 10).Value = 550
$ws.Cells.Item(18, 11).Value = 562.5
$ws.Cells.Item(18, 12).Value = 550
$ws.Cells.Item(18, 13).Value = -278.5
$ws.Cells.Item(18, 14).Value = -1118
$ws.Cells.Item(32, 8).Value = 437.36365
$ws.Cells.Item(32, 9).Value = 496.33334
$ws.Cells.Item(32, 10).Value = 172
$ws.Cells.Item(32, 11).Value = 496.33334
$ws.Cells.Item(32, 12).Value = 172
$ws.Cells.Item(32, 13).Value = -170.33334
$ws.Cells.Item(32, 14).Value = -824
$ws.Cells.Item(33, 8).Value = 625.7353000000001
$ws.Cells.Item(33, 9).Value = 622.65625
$ws.Cells.Item(33, 11).Value = 622.65625
$ws.Cells.Item(33, 13).Value = -393.65625
$ws.Cells.Item(34, 8).Value = 19552.455
$ws.Cells.Item(34, 9).Value = 1816.4
$ws.Cells.Item(34, 10).Value = 34332.5
$ws.Cells.Item(34, 11).Value = 1816.4
$ws.Cells.Item(34, 12).Value = 34332.5
$ws.Cells.Item(34, 13).Value = -1613.4
$ws.Cells.Item(34, 14).Value = -34738.5
$ws.Cells.Item(36, 8).Value = 19552.455
$ws.Cells.Item(36, 9).Value = 1816.4
$ws.Cells.Item(36, 10).Value = 34332.5
$ws.Cells.Item(36, 11).Value = 1816.4
$ws.Cells.Item(36, 12).Value = 34332.5
$ws.Cells.Item(36, 13).Value = -1101.4
$ws.Cells.Item(36, 14).Value = -35762.5
$ws.Cells.Item(40, 8).Value = 2118
$ws.Cells.Item(40, 9).Value = 1667.1428
$ws.Cells.Item(40, 10).Value = 2512.5
$ws.Cells.Item(40, 11).Value = 1667.1428
$ws.Cells.Item(40, 12).Value = 2512.5
$ws.Cells.Item(40, 13).Value = -1492.1428
$ws.Cells.Item(40, 14).Value = -2862.5
$ws.Cells.Item(43, 8).Value = 1074.5625
$ws.Cells.Item(43, 9).Value = 1258.5
$ws.Cells.Item(43, 10).Value = 964.2
$ws.Cells.Item(43, 11).Value = 1258.5
$ws.Cells.Item(43, 12).Value = 964.2
$ws.Cells.Item(43, 13).Value = -1189.5
$ws.Cells.Item(43, 14).Value = -1102.2
$ws.Cells.Item(101, 8).Value = 1799.3334
$ws.Cells.Item(101, 9).Value = 299
$ws.Cells.Item(101, 10).Value = 4800
$ws.Cells.Item(101, 11).Value = 897
$ws.Cells.Item(101, 12).Value = 14400
$ws.Cells.Item(101, 13).Value = 725
$ws.Cells.Item(101, 14).Value = -17644
$ws.Cells.Item(112, 8).Value = 1464.6923
$ws.Cells.Item(112, 10).Value = 1604.1
$ws.Cells.Item(112, 12).Value = 4812.299999999999
$ws.Cells.Item(112, 14).Value = -7028.299999999999
$ws.Cells.Item(129, 8).Value = 833.1
$ws.Cells.Item(129, 10).Value = 1103.35
$ws.Cells.Item(129, 12).Value = 3310.05
$ws.Cells.Item(129, 14).Value = -13310.05

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 9852.663
$ws.Cells.Item(32, 9).Value = 3393.4143
$ws.Cells.Item(32, 11).Value = 3393.4143
$ws.Cells.Item(32, 13).Value = -3106.4143
$ws.Cells.Item(61, 8).Value = 1277.9333
$ws.Cells.Item(61, 9).Value = 1226.6666
$ws.Cells.Item(61, 10).Value = 1483
$ws.Cells.Item(61, 11).Value = 1226.6666
$ws.Cells.Item(61, 12).Value = 1483
$ws.Cells.Item(61, 13).Value = -1014.6666
$ws.Cells.Item(61, 14).Value = -1907
$ws.Cells.Item(122, 8).Value = 1698.8823
$ws.Cells.Item(122, 9).Value = 911
$ws.Cells.Item(122, 10).Value = 3887.4443
$ws.Cells.Item(122, 11).Value = 2733
$ws.Cells.Item(122, 12).Value = 11662.3329
$ws.Cells.Item(122, 13).Value = -283
$ws.Cells.Item(122, 14).Value = -16562.3329
$ws.Cells.Item(132, 8).Value = 1380.2373
$ws.Cells.Item(132, 9).Value = 1111.5
$ws.Cells.Item(132, 10).Value = 2168.5334
$ws.Cells.Item(132, 11).Value = 3334.5
$ws.Cells.Item(132, 12).Value = 6505.600199999999
$ws.Cells.Item(132, 13).Value = -804.5
$ws.Cells.Item(132, 14).Value = -11565.6002
$ws.Cells.Item(136, 8).Value = 1277.9333
$ws.Cells.Item(136, 9).Value = 1226.6666
$ws.Cells.Item(136, 10).Value = 1483
$ws.Cells.Item(136, 11).Value = 3679.9998
$ws.Cells.Item(136, 12).Value = 4449
$ws.Cells.Item(136, 13).Value = -1129.9998
$ws.Cells.Item(136, 14).Value = -9549

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 31.333334
$ws.Cells.Item(22, 9).Value = 26.222221
$ws.Cells.Item(22, 10).Value = 46.666668
$ws.Cells.Item(22, 11).Value = 26.222221
$ws.Cells.Item(22, 12).Value = 46.666668
$ws.Cells.Item(22, 13).Value = 146.777779
$ws.Cells.Item(22, 14).Value = -392.666668
$ws.Cells.Item(80, 8).Value = 47.882355
$ws.Cells.Item(80, 10).Value = 53.57143
$ws.Cells.Item(80, 12).Value = 53.57143
$ws.Cells.Item(80, 14).Value = -2049.57143
$ws.Cells.Item(83, 8).Value = 47.882355
$ws.Cells.Item(83, 10).Value = 53.57143
$ws.Cells.Item(83, 12).Value = 267.85715
$ws.Cells.Item(83, 14).Value = -10251.85715
$ws.Cells.Item(99, 8).Value = 1914.2858
$ws.Cells.Item(99, 9).Value = 1950
$ws.Cells.Item(99, 10).Value = 1866.6666
$ws.Cells.Item(99, 11).Value = 1950
$ws.Cells.Item(99, 12).Value = 1866.6666
$ws.Cells.Item(99, 13).Value = -452
$ws.Cells.Item(99, 14).Value = -4862.6666
$ws.Cells.Item(134, 8).Value = 1144.9056
$ws.Cells.Item(134, 9).Value = 932.2449
$ws.Cells.Item(134, 10).Value = 3750
$ws.Cells.Item(134, 11).Value = 2796.7347
$ws.Cells.Item(134, 12).Value = 11250
$ws.Cells.Item(134, 13).Value = -261.7347
$ws.Cells.Item(134, 14).Value = -16320

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 57.714287
$ws.Cells.Item(7, 9).Value = 32.333332
$ws.Cells.Item(7, 10).Value = 76.75
$ws.Cells.Item(7, 11).Value = 32.333332
$ws.Cells.Item(7, 12).Value = 76.75
$ws.Cells.Item(7, 13).Value = 80.666668
$ws.Cells.Item(7, 14).Value = -302.75
$ws.Cells.Item(22, 8).Value = 442
$ws.Cells.Item(22, 9).Value = 477.5
$ws.Cells.Item(22, 10).Value = 300
$ws.Cells.Item(22, 11).Value = 477.5
$ws.Cells.Item(22, 12).Value = 300
$ws.Cells.Item(22, 13).Value = -127.5
$ws.Cells.Item(22, 14).Value = -1000

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 1786.6666
$ws.Cells.Item(113, 9).Value = 0
$ws.Cells.Item(113, 10).Value = 1786.6666
$ws.Cells.Item(113, 11).Value = 0
$ws.Cells.Item(113, 12).Value = 5359.9998
$ws.Cells.Item(113, 13).ClearContents()
$ws.Cells.Item(113, 14).Value = -9699.9998
$ws.Cells.Item(131, 8).Value = 1035.9108
$ws.Cells.Item(131, 9).Value = 464.75
$ws.Cells.Item(131, 10).Value = 1131.1041
$ws.Cells.Item(131, 11).Value = 1394.25
$ws.Cells.Item(131, 12).Value = 3393.3123
$ws.Cells.Item(131, 13).Value = 3645.75
$ws.Cells.Item(131, 14).Value = -13473.3123

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 108.42857
$ws.Cells.Item(2, 9).Value = 103.4
$ws.Cells.Item(2, 10).Value = 113
$ws.Cells.Item(2, 11).Value = 103.4
$ws.Cells.Item(2, 12).Value = 113
$ws.Cells.Item(2, 13).Value = 9.599999999999994
$ws.Cells.Item(2, 14).Value = -339
$ws.Cells.Item(22, 8).Value = 2000
$ws.Cells.Item(22, 10).Value = 2000
$ws.Cells.Item(22, 12).Value = 2000
$ws.Cells.Item(22, 14).Value = -3058

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1667204.6
$ws.Cells.Item(22, 9).Value = 4762102.5
$ws.Cells.Item(22, 10).Value = 721.2308
$ws.Cells.Item(22, 11).Value = 4762102.5
$ws.Cells.Item(22, 12).Value = 721.2308
$ws.Cells.Item(22, 13).Value = -4761807.5
$ws.Cells.Item(22, 14).Value = -1311.2308
$ws.Cells.Item(27, 8).Value = 1667204.6
$ws.Cells.Item(27, 9).Value = 4762102.5
$ws.Cells.Item(27, 10).Value = 721.2308
$ws.Cells.Item(27, 11).Value = 4762102.5
$ws.Cells.Item(27, 12).Value = 721.2308
$ws.Cells.Item(27, 13).Value = -4761995.5
$ws.Cells.Item(27, 14).Value = -935.2308

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 1136.0588
$ws.Cells.Item(136, 9).Value = 639.8276
$ws.Cells.Item(136, 10).Value = 1790.1818
$ws.Cells.Item(136, 11).Value = 1919.4828
$ws.Cells.Item(136, 12).Value = 5370.5454
$ws.Cells.Item(136, 13).Value = 630.5172000000002
$ws.Cells.Item(136, 14).Value = -10470.5454

Write-Output "Applied all updates"